$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for rows 2-8
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = -4
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 1
